# Implemented getting kafka relations.
# Reorders the rows of the "classFields" sheet (Order class field metadata)
# so that Field Name / Field Type pairs appear in the new order shown by
# the diff, while Class Name / Field Modifier stay the same for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# New order (rows 2..8), each entry: Field Name -> Field Type
$rows = @(
    @{ Name = "source";       Type = "java.lang.String" },
    @{ Name = "productId";    Type = "java.lang.Long" },
    @{ Name = "price";        Type = "int" },
    @{ Name = "status";       Type = "java.lang.String" },
    @{ Name = "id";           Type = "java.lang.Long" },
    @{ Name = "productCount"; Type = "int" },
    @{ Name = "customerId";   Type = "java.lang.Long" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 4).Value = $row.Type
    $r = $r + 1
}
